# Apply the "Updated cryptos list" refresh: new prices/volumes, and a few
# rows whose rank order swapped (name/link/price/volume move together).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.840.14"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.584.20"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "3.049.01"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "62.788.65"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000146"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "2.593.45"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("D24").Value = "2.715.49"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0818"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "466.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.402"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "158.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.629"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.23%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0540"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0965"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
